$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.387.28"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.23"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7135"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.42"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08030"
$ws.Range("E8").Value = "  +3.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3130"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.28"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08364"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.916.46"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7212"
$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.248"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.67"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.292"
$ws.Range("E16").Value = "  +4.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008478"
$ws.Range("E17").Value = "  +2.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.396.01"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.31"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.26"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.130.97"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.854"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"

$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.65"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.075"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("E30").Value = "  +0.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.343"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.211"
$ws.Range("E32").Value = "  -5.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05379"
$ws.Range("E33").Value = "  +2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.954"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7501"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  +0.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.285.10"
$ws.Range("E39").Value = "  +9.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.749"
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.575"
$ws.Range("E41").Value = "  +2.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.46"
$ws.Range("E42").Value = "  +4.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.78"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8932"
$ws.Range("E44").Value = "  +0.80%  "

$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("E46").Value = "  +6.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.029.57"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.809"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5217"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.503"
$ws.Range("E50").Value = "  +1.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4396"
$ws.Range("E51").Value = "  +2.12%  "
